$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")

# Insert a new row at position 7 (pushes the existing row 7 (blank) and
# everything below it down by one row). Excel carries the formatting of the
# row above down into the freshly inserted row, matching the bold "field
# name" style in column A and the wrapped "field value" style in column B.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with the business therapeutic areas entry.
$ws.Cells.Item(7, 1).Value = "businessTherapeuticAreas"
$ws.Cells.Item(7, 2).Value = "SPONSOR: PHARMA=Pharma Division"

# Explicit (custom) row height, same as the sheet default.
$ws.Rows.Item(7).RowHeight = 16

# Update the recorded selection to match the post-edit state.
$ws.Range("B14").Select()
